# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Estado de Cuenta" worker / period / value rows (rows 16-25).
# Each worker now lists periods 2312 -> 2308 (most recent first) consecutively,
# replacing the previous interleaved arrangement.
$data = @(
    @{row=16; id="37617413";   name="YURISAN PATIÑO BOHORQUEZ"; periodo="2312"; valor=25333},
    @{row=17; id="37617413";   name="YURISAN PATIÑO BOHORQUEZ"; periodo="2311"; valor=46400},
    @{row=18; id="37617413";   name="YURISAN PATIÑO BOHORQUEZ"; periodo="2310"; valor=46400},
    @{row=19; id="37617413";   name="YURISAN PATIÑO BOHORQUEZ"; periodo="2309"; valor=46400},
    @{row=20; id="37617413";   name="YURISAN PATIÑO BOHORQUEZ"; periodo="2308"; valor=46400},
    @{row=21; id="1050544690"; name="WILDER SANJUAN SERRANO";   periodo="2312"; valor=25333},
    @{row=22; id="1050544690"; name="WILDER SANJUAN SERRANO";   periodo="2311"; valor=46400},
    @{row=23; id="1050544690"; name="WILDER SANJUAN SERRANO";   periodo="2310"; valor=46400},
    @{row=24; id="1050544690"; name="WILDER SANJUAN SERRANO";   periodo="2309"; valor=46400},
    @{row=25; id="1050544690"; name="WILDER SANJUAN SERRANO";   periodo="2308"; valor=46400}
)

foreach ($d in $data) {
    $r = $d.row
    $ws.Range("C$r").Value = $d.id
    $ws.Range("D$r").Value = $d.name
    $ws.Range("E$r").Value = $d.periodo
    $ws.Range("F$r").Value = $d.valor
}
